$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header in C1 (same style as B1)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("A2").Value = "1810 Renta variable"
$ws.Range("B2").Value = 88171.49000000001
$ws.Range("C2").Value = 87857.10000000001
$ws.Range("A3").Value = "1822 Raices Valores Negociables"
$ws.Range("B3").Value = 338229.12
$ws.Range("C3").Value = 355152.98
$ws.Range("A4").Value = "Adcap IOL Acciones Argentina"
$ws.Range("B4").Value = 40275.09
$ws.Range("C4").Value = 40345.22
$ws.Range("A5").Value = "Allaria Acciones"
$ws.Range("B5").Value = 61513.47
$ws.Range("C5").Value = 61431.38
$ws.Range("A6").Value = "Alpha Acciones"
$ws.Range("B6").Value = 160594.77
$ws.Range("C6").Value = 160738.57
$ws.Range("A7").Value = "Alpha Latam"
$ws.Range("B7").Value = 59.04
$ws.Range("C7").Value = 56.35
$ws.Range("A8").Value = "Alpha Mega"
$ws.Range("B8").Value = 147814.31
$ws.Range("C8").Value = 148451.41
$ws.Range("A9").Value = "Alpha renta balan global"
$ws.Range("B9").Value = 279887.82
$ws.Range("C9").Value = 279976.7
$ws.Range("A10").Value = "Argenfunds"
$ws.Range("B10").Value = 12884.17
$ws.Range("C10").Value = 12891.26
$ws.Range("A11").Value = "Arpenta acciones"
$ws.Range("B11").Value = 2995.62
$ws.Range("C11").Value = 2993.91
$ws.Range("A12").Value = "Arpenta ex Mercosur"
$ws.Range("B12").Value = 1979.47
$ws.Range("C12").Value = 1978.28
$ws.Range("A13").Value = "Balanz"
$ws.Range("B13").Value = 247792.75
$ws.Range("C13").Value = 160513.36
$ws.Range("A14").Value = "Compass Crecimiento"
$ws.Range("B14").Value = 1145282.22
$ws.Range("C14").Value = 1130374.27
$ws.Range("A15").Value = "Compass Crecimiento II"
$ws.Range("B15").Value = 17910.03
$ws.Range("C15").Value = 17914.94
$ws.Range("A16").Value = "Consultatio Acciones Argentina"
$ws.Range("B16").Value = 1141958.47
$ws.Range("C16").Value = 1142049.67
$ws.Range("A17").Value = "Consultatio Renta Variable"
$ws.Range("B17").Value = 420944.4
$ws.Range("C17").Value = 420777.72
$ws.Range("A18").Value = "Delta Acciones"
$ws.Range("B18").Value = 75202.55
$ws.Range("C18").Value = 75304
$ws.Range("A19").Value = "Delta Internacional"
$ws.Range("B19").Value = 2504.4
$ws.Range("C19").Value = 2498.91
$ws.Range("A20").Value = "Delta Latinoamerica"
$ws.Range("B20").Value = 8194.18
$ws.Range("C20").Value = 8186.09
$ws.Range("A21").Value = "Delta Recursos Naturales"
$ws.Range("B21").Value = 34941.85
$ws.Range("C21").Value = 34986.61
$ws.Range("A22").Value = "Delta Select"
$ws.Range("B22").Value = 442216.7
$ws.Range("C22").Value = 442441.07
$ws.Range("A23").Value = "Delta gestion V"
$ws.Range("B23").Value = 106360.7
$ws.Range("C23").Value = 106795.6
$ws.Range("A24").Value = "FBA Acciones Argentinas"
$ws.Range("B24").Value = 234004.44
$ws.Range("C24").Value = 239614.19
$ws.Range("A25").Value = "FBA Calificado"
$ws.Range("B25").Value = 230700.79
$ws.Range("C25").Value = 234272.04
$ws.Range("A26").Value = "Fima Acciones"
$ws.Range("B26").Value = 249184.42
$ws.Range("C26").Value = 257751.14
$ws.Range("A27").Value = "Fima PB Acciones"
$ws.Range("B27").Value = 206464.83
$ws.Range("C27").Value = 216929.81
$ws.Range("A28").Value = "Galileo Acciones"
$ws.Range("B28").Value = 2073702.77
$ws.Range("C28").Value = 2149584.31
$ws.Range("A29").Value = "Goal Acciones Argentinas"
$ws.Range("B29").Value = 40034.88
$ws.Range("C29").Value = 39995.45
$ws.Range("A30").Value = "Goal acciones plus"
$ws.Range("B30").Value = 5291.3
$ws.Range("C30").Value = 5304.06
$ws.Range("A31").Value = "HF Acciones Argentinas"
$ws.Range("B31").Value = 157899.35
$ws.Range("C31").Value = 158126.71
$ws.Range("A32").Value = "HF Acciones Lideres"
$ws.Range("B32").Value = 300647.63
$ws.Range("C32").Value = 295801.43
$ws.Range("A33").Value = "IAM Renta Variable"
$ws.Range("B33").Value = 41744.48
$ws.Range("C33").Value = 43603.61
$ws.Range("A34").Value = "IEB Value"
$ws.Range("B34").Value = 3896.97
$ws.Range("C34").Value = 3897.54
$ws.Range("A35").Value = "Lombardi"
$ws.Range("B35").Value = 37560.76
$ws.Range("C35").Value = 37467.97
$ws.Range("A36").Value = "MAF"
$ws.Range("B36").Value = 79675.07000000001
$ws.Range("C36").Value = 79731.38
$ws.Range("A37").Value = "Megainver"
$ws.Range("B37").Value = 35003.74
$ws.Range("C37").Value = 35010.68
$ws.Range("A38").Value = "Pellegrini Acciones"
$ws.Range("B38").Value = 85115.45
$ws.Range("C38").Value = 85065.67999999999
$ws.Range("A39").Value = "Pionero Acciones"
$ws.Range("B39").Value = 129626.95
$ws.Range("C39").Value = 129482.29
$ws.Range("A40").Value = "Premier Renta Variable"
$ws.Range("B40").Value = 69373
$ws.Range("C40").Value = 69275.31
$ws.Range("A41").Value = "Quinquela Acciones"
$ws.Range("B41").Value = 117022.34
$ws.Range("C41").Value = 116859.7
$ws.Range("A42").Value = "Rofex 20 Renta Variable"
$ws.Range("B42").Value = 82812.8
$ws.Range("C42").Value = 82783.25
$ws.Range("A43").Value = "SBS Acciones Argentina"
$ws.Range("B43").Value = 576713.58
$ws.Range("C43").Value = 576814.16
$ws.Range("A44").Value = "Schroeder RV"
$ws.Range("B44").Value = 1211874.1
$ws.Range("C44").Value = 1212773.08
$ws.Range("A45").Value = "Supefondo RV"
$ws.Range("B45").Value = 615981.92
$ws.Range("C45").Value = 617718.1
$ws.Range("A46").Value = "Superfondo "
$ws.Range("B46").Value = 16264.35
$ws.Range("C46").Value = 17100.51
$ws.Range("A47").Value = "Toronto Trust Multimercado"
$ws.Range("B47").Value = 60630.18
$ws.Range("C47").Value = 60598.85
$ws.Range("A48").Value = "Toronto trust Argy"
$ws.Range("B48").Value = 40119.38
$ws.Range("C48").Value = 67655.60000000001
$ws.Range("A49").Value = "avg"
$ws.Range("B49").Value = 244235.28
$ws.Range("C49").Value = 245253.88
$ws.Range("A50").Value = "total"
$ws.Range("B50").Value = 11479058.1
$ws.Range("C50").Value = 11526932.25
